# Add a "Total" row (13) under the product table and a SUM formula,
# then apply a two-color-scale conditional format over A2:A11
# (added twice, matching the duplicated cfRule in the target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: label + total formula
$ws.Range("A13").Value = "Total"
$ws.Range("B13").Formula = "=SUM(A2:A11)"

# Conditional formatting: color scale over A2:A11 (duplicated rule)
$range = $ws.Range("A2:A11")
$cf1 = $range.FormatConditions.AddColorScale(2)
$cf1.ColorScaleCriteria.Item(1).FormatColor.Color = 16777184
$cf1.ColorScaleCriteria.Item(2).FormatColor.Color = 8421376

$cf2 = $range.FormatConditions.AddColorScale(2)
$cf2.ColorScaleCriteria.Item(1).FormatColor.Color = 16777184
$cf2.ColorScaleCriteria.Item(2).FormatColor.Color = 8421376
